$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991109864/1750183612496/A5MTrC4t3W-camera-video-segment-driver-1750183614996.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSM6UWVVNW%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T200004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKn%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDXddvqJrGWvCXWEhjFTvw8Cfh2QH66%2Bj1UD12%2BPcf5zQIgXH%2BHEaeQfjTfN7twEyNKi6%2FyiSLD6aFG9sZebDHxlskq5gMIkv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDG3hRqqGNpvelMiIzSq6A83pocEq75MbkMG753ojvtDqBTYDaoy%2B3wh8UD%2B0Q30lLtedzIqyNQIT3S%2FnVyXBQEpDKuc04hgP9hyw%2BNkMTH0stGAMpNlKd3hM0aWY6dzyb2KWIfSBB8gELf9xjgEB7BMlsqux8sMw4TUeC1tdefA44XBJGKjwXH07KdBnAubzoFkeKAmHKF8cgW%2Fx6AqXTVS8VZmGkZjVSYaPjFeGBgg%2BnwCB%2BNCRAHsDRVvOSUaKvMsTOvXeP88zc5wfbe3wAuGPK1Sng4yf33ji9jF%2FbQPZNZCTNYCLGBgEnsGsYv0J6u58YQrgINqn2vclfjxlxmIaAKpRpjEQyv57X4RRSNZPnT%2BSryliz0ETl5j7wkhBkZDOLlXINQWtgOjZkQ0deUvAJZsXbC0M9uTMqdqNf%2BLe67lUmsv5JFeLrwGjBUX48s%2FU%2BNn4GOZ4zLaERZL%2FJ5Be4oKZIUKycE%2FOouV0xH8Bui6PDJkOWsVmyV2vScWsoZhWLFOx%2FtlETcwp76bcAQ0rImY%2FkJDMt3FnF38%2FXtjhd40JxvURMyQtMahxiB1ypz6d0dzRIzwtAUOOIsrPf%2FEf7zTfbbFmFsowht3LwgY6pQGgNv%2F7BRaALSce2r%2B1CYeXmi5unBmfbNNnYg3KyFwBOELPRygwGpFt%2FHg70vMlt4yVVOQUegpVKMBZTGYstMw%2Bo6oCB6cb5JQyCC5%2BFIqpebNNIre6a9UTJw6LXtqIe3n4OnpEPG9re97um9KQSXczgxTQi7U3UMD9lbqi7AvF2OBqExYWB4wEi1GWMTW7uYgBuRk6dNMvumkMg8i9ohxFm6GtK14%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2004%3A00%3A04%20GMT&X-Amz-Signature=bd97b82f5b60c7fc9d2e14440f7fd89200dad783304b37401bdab3375762a5e9"

$ws.Range("L2").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1750183612496/TAdN8NokD3-camera-video-segment-1750183614996.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSM6UWVVNW%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T200004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKn%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDXddvqJrGWvCXWEhjFTvw8Cfh2QH66%2Bj1UD12%2BPcf5zQIgXH%2BHEaeQfjTfN7twEyNKi6%2FyiSLD6aFG9sZebDHxlskq5gMIkv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDG3hRqqGNpvelMiIzSq6A83pocEq75MbkMG753ojvtDqBTYDaoy%2B3wh8UD%2B0Q30lLtedzIqyNQIT3S%2FnVyXBQEpDKuc04hgP9hyw%2BNkMTH0stGAMpNlKd3hM0aWY6dzyb2KWIfSBB8gELf9xjgEB7BMlsqux8sMw4TUeC1tdefA44XBJGKjwXH07KdBnAubzoFkeKAmHKF8cgW%2Fx6AqXTVS8VZmGkZjVSYaPjFeGBgg%2BnwCB%2BNCRAHsDRVvOSUaKvMsTOvXeP88zc5wfbe3wAuGPK1Sng4yf33ji9jF%2FbQPZNZCTNYCLGBgEnsGsYv0J6u58YQrgINqn2vclfjxlxmIaAKpRpjEQyv57X4RRSNZPnT%2BSryliz0ETl5j7wkhBkZDOLlXINQWtgOjZkQ0deUvAJZsXbC0M9uTMqdqNf%2BLe67lUmsv5JFeLrwGjBUX48s%2FU%2BNn4GOZ4zLaERZL%2FJ5Be4oKZIUKycE%2FOouV0xH8Bui6PDJkOWsVmyV2vScWsoZhWLFOx%2FtlETcwp76bcAQ0rImY%2FkJDMt3FnF38%2FXtjhd40JxvURMyQtMahxiB1ypz6d0dzRIzwtAUOOIsrPf%2FEf7zTfbbFmFsowht3LwgY6pQGgNv%2F7BRaALSce2r%2B1CYeXmi5unBmfbNNnYg3KyFwBOELPRygwGpFt%2FHg70vMlt4yVVOQUegpVKMBZTGYstMw%2Bo6oCB6cb5JQyCC5%2BFIqpebNNIre6a9UTJw6LXtqIe3n4OnpEPG9re97um9KQSXczgxTQi7U3UMD9lbqi7AvF2OBqExYWB4wEi1GWMTW7uYgBuRk6dNMvumkMg8i9ohxFm6GtK14%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2004%3A00%3A04%20GMT&X-Amz-Signature=24dcea7ab6c468f59496521c0a1a8546d7bab790671874cddfff6455078b4f06"

$ws.Range("K3").Value = "https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991152736/1750182961683/iBjUYEdqdu-camera-video-segment-driver-1750182964183.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSM6UWVVNW%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T200004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKn%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDXddvqJrGWvCXWEhjFTvw8Cfh2QH66%2Bj1UD12%2BPcf5zQIgXH%2BHEaeQfjTfN7twEyNKi6%2FyiSLD6aFG9sZebDHxlskq5gMIkv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDG3hRqqGNpvelMiIzSq6A83pocEq75MbkMG753ojvtDqBTYDaoy%2B3wh8UD%2B0Q30lLtedzIqyNQIT3S%2FnVyXBQEpDKuc04hgP9hyw%2BNkMTH0stGAMpNlKd3hM0aWY6dzyb2KWIfSBB8gELf9xjgEB7BMlsqux8sMw4TUeC1tdefA44XBJGKjwXH07KdBnAubzoFkeKAmHKF8cgW%2Fx6AqXTVS8VZmGkZjVSYaPjFeGBgg%2BnwCB%2BNCRAHsDRVvOSUaKvMsTOvXeP88zc5wfbe3wAuGPK1Sng4yf33ji9jF%2FbQPZNZCTNYCLGBgEnsGsYv0J6u58YQrgINqn2vclfjxlxmIaAKpRpjEQyv57X4RRSNZPnT%2BSryliz0ETl5j7wkhBkZDOLlXINQWtgOjZkQ0deUvAJZsXbC0M9uTMqdqNf%2BLe67lUmsv5JFeLrwGjBUX48s%2FU%2BNn4GOZ4zLaERZL%2FJ5Be4oKZIUKycE%2FOouV0xH8Bui6PDJkOWsVmyV2vScWsoZhWLFOx%2FtlETcwp76bcAQ0rImY%2FkJDMt3FnF38%2FXtjhd40JxvURMyQtMahxiB1ypz6d0dzRIzwtAUOOIsrPf%2FEf7zTfbbFmFsowht3LwgY6pQGgNv%2F7BRaALSce2r%2B1CYeXmi5unBmfbNNnYg3KyFwBOELPRygwGpFt%2FHg70vMlt4yVVOQUegpVKMBZTGYstMw%2Bo6oCB6cb5JQyCC5%2BFIqpebNNIre6a9UTJw6LXtqIe3n4OnpEPG9re97um9KQSXczgxTQi7U3UMD9lbqi7AvF2OBqExYWB4wEi1GWMTW7uYgBuRk6dNMvumkMg8i9ohxFm6GtK14%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2004%3A00%3A04%20GMT&X-Amz-Signature=ba960bde0bf46ff01b88360206d16b2093b471e10334f100d5d88663d341e596"

$ws.Range("L3").Value = "https://s3.samsara.com/samsara-cvdata/4006124/281474991152736/1750182961683/Foph7ZSROm-camera-video-segment-1750182964183.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSM6UWVVNW%2F20250618%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250618T200004Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEKn%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQDXddvqJrGWvCXWEhjFTvw8Cfh2QH66%2Bj1UD12%2BPcf5zQIgXH%2BHEaeQfjTfN7twEyNKi6%2FyiSLD6aFG9sZebDHxlskq5gMIkv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDG3hRqqGNpvelMiIzSq6A83pocEq75MbkMG753ojvtDqBTYDaoy%2B3wh8UD%2B0Q30lLtedzIqyNQIT3S%2FnVyXBQEpDKuc04hgP9hyw%2BNkMTH0stGAMpNlKd3hM0aWY6dzyb2KWIfSBB8gELf9xjgEB7BMlsqux8sMw4TUeC1tdefA44XBJGKjwXH07KdBnAubzoFkeKAmHKF8cgW%2Fx6AqXTVS8VZmGkZjVSYaPjFeGBgg%2BnwCB%2BNCRAHsDRVvOSUaKvMsTOvXeP88zc5wfbe3wAuGPK1Sng4yf33ji9jF%2FbQPZNZCTNYCLGBgEnsGsYv0J6u58YQrgINqn2vclfjxlxmIaAKpRpjEQyv57X4RRSNZPnT%2BSryliz0ETl5j7wkhBkZDOLlXINQWtgOjZkQ0deUvAJZsXbC0M9uTMqdqNf%2BLe67lUmsv5JFeLrwGjBUX48s%2FU%2BNn4GOZ4zLaERZL%2FJ5Be4oKZIUKycE%2FOouV0xH8Bui6PDJkOWsVmyV2vScWsoZhWLFOx%2FtlETcwp76bcAQ0rImY%2FkJDMt3FnF38%2FXtjhd40JxvURMyQtMahxiB1ypz6d0dzRIzwtAUOOIsrPf%2FEf7zTfbbFmFsowht3LwgY6pQGgNv%2F7BRaALSce2r%2B1CYeXmi5unBmfbNNnYg3KyFwBOELPRygwGpFt%2FHg70vMlt4yVVOQUegpVKMBZTGYstMw%2Bo6oCB6cb5JQyCC5%2BFIqpebNNIre6a9UTJw6LXtqIe3n4OnpEPG9re97um9KQSXczgxTQi7U3UMD9lbqi7AvF2OBqExYWB4wEi1GWMTW7uYgBuRk6dNMvumkMg8i9ohxFm6GtK14%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2019%20Jun%202025%2004%3A00%3A04%20GMT&X-Amz-Signature=f0605978ce4381f568277d578424439344ae21ad32bc35a5198c0dfc8426b3b5"
